$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.422.19'
$ws.Range("E2").Value = '  +0.85%  '

$ws.Range("D3").Value = '1.875.32'
$ws.Range("E3").Value = '  -0.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.98%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.015'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5136'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3928'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08322'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.120'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.95'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.266'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("D13").Value = '1.873.22'
$ws.Range("E13").Value = '  +0.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.246'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.014'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001106'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06726'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.982'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.24%  '

$ws.Range("D23").Value = '28.464.02'
$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.257'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.25%  '

$ws.Range("D26").Value = '2.087.45'
$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.414'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1058'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.041'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.869'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.638'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02452'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06531'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.129'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2186'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.258'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6468'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.187'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.001'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.81%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6038'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.706'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.285'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.005'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.210'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06889'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '

